$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update masthead rich text (issue number + date range) ---
$ws.Range("A8").Value = "Volume 31   Number  10"
$ws.Range("C9").Value = "Report Covering the Week  3/4/2024  Through  3/10/2024"

# --- Cells flipping from the "0" placeholder (text) to a real numeric value ---
# Copy number format from a neighboring numeric cell, then write the value.
$ws.Range("G15").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Value = 1
$ws.Range("G15").Copy() | Out-Null
$ws.Range("F15").PasteSpecial(-4122) | Out-Null
$ws.Range("F15").Value = 1
$ws.Range("J15").Copy() | Out-Null
$ws.Range("I15").PasteSpecial(-4122) | Out-Null
$ws.Range("I15").Value = 1
$ws.Range("D20").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$ws.Range("C20").Value = 1

# --- Cells flipping from a real numeric value back to the "0" placeholder (text) ---
# Copy value+format from the neighboring placeholder cell, then re-apply its format
# (PasteSpecial Values-then-Formats keeps the shared "0" text and the dash-style format).
$ws.Range("G29").Copy() | Out-Null
$ws.Range("F29").PasteSpecial(-4163) | Out-Null
$ws.Range("G29").Copy() | Out-Null
$ws.Range("F29").PasteSpecial(-4122) | Out-Null
$ws.Range("G30").Copy() | Out-Null
$ws.Range("F30").PasteSpecial(-4163) | Out-Null
$ws.Range("G30").Copy() | Out-Null
$ws.Range("F30").PasteSpecial(-4122) | Out-Null
$ws.Range("G31").Copy() | Out-Null
$ws.Range("F31").PasteSpecial(-4163) | Out-Null
$ws.Range("G31").Copy() | Out-Null
$ws.Range("F31").PasteSpecial(-4122) | Out-Null

# --- Plain numeric value updates ---
$ws.Range("N14").Value = -50
$ws.Range("H15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = -80
$ws.Range("M15").Value = -50
$ws.Range("N15").Value = -66.666666666666
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 15
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 32
$ws.Range("J16").Value = 35
$ws.Range("K16").Value = -8.571428571428
$ws.Range("L16").Value = -17.948717948717
$ws.Range("N16").Value = -82.320441988950
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 8
$ws.Range("H17").Value = -46.666666666666
$ws.Range("I17").Value = 28
$ws.Range("J17").Value = 41
$ws.Range("K17").Value = -31.707317073170
$ws.Range("L17").Value = -34.883720930232
$ws.Range("M17").Value = 64.705882352941
$ws.Range("N17").Value = -36.363636363636
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 400
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 46
$ws.Range("J18").Value = 41
$ws.Range("K18").Value = 12.195121951219
$ws.Range("L18").Value = 4.545454545454
$ws.Range("M18").Value = -34.285714285714
$ws.Range("N18").Value = -86.350148367952
$ws.Range("C19").Value = 20
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = 5.263157894736
$ws.Range("F19").Value = 56
$ws.Range("G19").Value = 77
$ws.Range("H19").Value = -27.272727272727
$ws.Range("I19").Value = 175
$ws.Range("J19").Value = 180
$ws.Range("K19").Value = -2.777777777777
$ws.Range("L19").Value = -0.568181818181
$ws.Range("M19").Value = -31.102362204724
$ws.Range("N19").Value = -64.574898785425
$ws.Range("E20").Value = 0
$ws.Range("I20").Value = 5
$ws.Range("J20").Value = 10
$ws.Range("K20").Value = -50
$ws.Range("L20").Value = -44.444444444444
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -97.959183673469
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = 28
$ws.Range("F21").Value = 100
$ws.Range("G21").Value = 128
$ws.Range("H21").Value = -21.875
$ws.Range("I21").Value = 288
$ws.Range("J21").Value = 308
$ws.Range("K21").Value = -6.493506493506
$ws.Range("L21").Value = -8.860759493670
$ws.Range("M21").Value = -24.210526315789
$ws.Range("N21").Value = -77.947932618683
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = -80
$ws.Range("J22").Value = 18
$ws.Range("K22").Value = -11.111111111111
$ws.Range("L22").Value = -30.434782608695
$ws.Range("C24").Value = 49
$ws.Range("D24").Value = 40
$ws.Range("E24").Value = 22.5
$ws.Range("F24").Value = 249
$ws.Range("G24").Value = 134
$ws.Range("H24").Value = 85.820895522388
$ws.Range("I24").Value = 580
$ws.Range("J24").Value = 365
$ws.Range("K24").Value = 58.904109589041
$ws.Range("L24").Value = 34.883720930232
$ws.Range("M24").Value = 80.685358255451
$ws.Range("C25").Value = 39
$ws.Range("D25").Value = 25
$ws.Range("E25").Value = 56
$ws.Range("F25").Value = 201
$ws.Range("G25").Value = 102
$ws.Range("H25").Value = 97.058823529411
$ws.Range("I25").Value = 490
$ws.Range("J25").Value = 264
$ws.Range("K25").Value = 85.606060606060
$ws.Range("L25").Value = 57.051282051282
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -9.090909090909
$ws.Range("F26").Value = 53
$ws.Range("G26").Value = 38
$ws.Range("H26").Value = 39.473684210526
$ws.Range("I26").Value = 103
$ws.Range("J26").Value = 81
$ws.Range("K26").Value = 27.160493827160
$ws.Range("L26").Value = 19.767441860465
$ws.Range("M26").Value = 53.731343283582
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 3
$ws.Range("K27").Value = -25
$ws.Range("L27").Value = -57.142857142857
$ws.Range("D28").Value = 7
$ws.Range("E28").Value = -85.714285714285
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 11
$ws.Range("H28").Value = -63.636363636363
$ws.Range("J28").Value = 20
$ws.Range("K28").Value = -25
$ws.Range("L28").Value = -11.764705882352
$ws.Range("N29").Value = -50
$ws.Range("N30").Value = -50
$ws.Range("I31").Value = 3
$ws.Range("K31").Value = 200
$ws.Range("L31").Value = -57.142857142857

Write-Host "edit complete"
